# Auto-generated Excel COM-interop edit script
# Implements: update report header date, extend Report sheet data through
# row 97 (5 new work-order rows), adjust two cell wrap styles on row 92,
# extend the print area, and move the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# ---------------------------------------------------------------------
# 1. Update the report title/header (shared string index 0) with the new
#    "製表日期" (report-generation date): 2025-09-17 -> 2025-09-18
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202509   (  製表日期:2025-09-18  )"

# ---------------------------------------------------------------------
# 2. Row 92 gains wrap-text formatting on P92 (still blank) and AC92
#    (existing "裝潢回裝完成" note).
# ---------------------------------------------------------------------
$ws.Range("P92").WrapText = $true
$ws.Range("AC92").WrapText = $true

# ---------------------------------------------------------------------
# 3. Seed rows 93-97 by copying the nearest same-banding template row so
#    the alternating row fill / borders / alignment carry over exactly,
#    then overwrite the values cell-by-cell.
#    Row 93 <- row 3 (odd/banded template)
#    Row 94 <- row 4 (even/plain template)
#    Row 95 <- row 3
#    Row 96 <- row 4
#    Row 97 <- row 3
# ---------------------------------------------------------------------
$ws.Range("A3:AK3").Copy($ws.Range("A93:AK93"))
$ws.Range("A4:AK4").Copy($ws.Range("A94:AK94"))
$ws.Range("A3:AK3").Copy($ws.Range("A95:AK95"))
$ws.Range("A4:AK4").Copy($ws.Range("A96:AK96"))
$ws.Range("A3:AK3").Copy($ws.Range("A97:AK97"))

$ws.Range("A93").Value = 91
$ws.Range("B93").Value = "維修"
$ws.Range("C93").Value = 2025092373
$ws.Range("D93").Value = "ED620114091801"
$ws.Range("E93").Value = "一般件"
$ws.Range("F93").Value = "D620"
$ws.Range("G93").Value = "三重福隆店"
$ws.Range("H93").Value = "新北市三重區"
$ws.Range("I93").Value = "2025-09-18 00:12:53"
$ws.Range("J93").Value = "星期四"
$ws.Range("K93").Value = "凌晨"
$ws.Range("L93").Value = "HLF2"
$ws.Range("M93").Value = "HL-CCD掃描器(TM)"
$ws.Range("N93").Value = "F201"
$ws.Range("O93").Value = "掃描無反應或感應不良"
$ws.Range("P93").Value = "門市反應TM1-CCD掃描器(HC76-TR)無電源反應，已確認後方線路無鬆脫，重啟TM仍異常..請台芝到店協助(掃描器無紅光，無法掃描)`n9/18 09:02 致電門市未接..吳"
$ws.Range("Q93").Value = "THILF0D620"
$ws.Range("R93").Value = "新北一"
$ws.Range("S93").Value = "吳宗鴻"
$ws.Range("T93").Value = 1
$ws.Range("U93").Value = "已完工"
$ws.Range("V93").Value = "2025-09-18 09:14:36"
$ws.Range("W93").Value = "2025-09-18 13:21:00"
$ws.Range("X93").Value = "2025-09-18 13:51:00"
$ws.Range("Y93").Value = "2025-09-19 13:14:00"
$ws.Range("Z93").Value = 0.5
$ws.Range("AA93").Value = ""
$ws.Range("AB93").Value = "到場處理"
$ws.Range("AC93").Value = "TX800的IO卡無反應，更換IO卡後測試正常"
$ws.Range("AD93").Value = ""
$ws.Range("AE93").Value = ""
$ws.Range("AF93").Value = ""
$ws.Range("AG93").Value = ""
$ws.Range("AH93").Value = ""
$ws.Range("AI93").Value = ""
$ws.Range("AJ93").Value = ""
$ws.Range("AK93").Value = "O"

$ws.Range("A94").Value = 92
$ws.Range("B94").Value = "維修"
$ws.Range("C94").Value = 2025092384
$ws.Range("D94").Value = "1D191114091801"
$ws.Range("E94").Value = "一般件"
$ws.Range("F94").Value = "D191"
$ws.Range("G94").Value = "三重興德店"
$ws.Range("H94").Value = "新北市三重區"
$ws.Range("I94").Value = "2025-09-18 10:03:23"
$ws.Range("J94").Value = "星期四"
$ws.Range("K94").Value = "上午"
$ws.Range("L94").Value = "HLF2"
$ws.Range("M94").Value = "HL-CCD掃描器(TM)"
$ws.Range("N94").Value = "F201"
$ws.Range("O94").Value = "掃描無反應或感應不良"
$ws.Range("P94").Value = "門市反應TM1 CCD掃描器(HC56II-TR)線路外層有脫落且刷讀條碼不太好刷，有執行校正仍異常....需請台芝到店協助"
$ws.Range("Q94").Value = "THILF0D191"
$ws.Range("R94").Value = "新北一"
$ws.Range("S94").Value = "吳宗鴻"
$ws.Range("T94").Value = 1
$ws.Range("U94").Value = "已完工"
$ws.Range("V94").Value = "2025-09-18 10:05:29"
$ws.Range("W94").Value = "2025-09-18 12:10:00"
$ws.Range("X94").Value = "2025-09-18 12:40:00"
$ws.Range("Y94").Value = "2025-09-19 14:05:00"
$ws.Range("Z94").Value = 0.5
$ws.Range("AA94").Value = ""
$ws.Range("AB94").Value = "到場處理"
$ws.Range("AC94").Value = "更換掃描槍`n換下8119008295`n換上8119012936"
$ws.Range("AD94").Value = ""
$ws.Range("AE94").Value = ""
$ws.Range("AF94").Value = ""
$ws.Range("AG94").Value = ""
$ws.Range("AH94").Value = ""
$ws.Range("AI94").Value = ""
$ws.Range("AJ94").Value = ""
$ws.Range("AK94").Value = "O"

$ws.Range("A95").Value = 93
$ws.Range("B95").Value = "服務"
$ws.Range("C95").Value = 2025092399
$ws.Range("D95").Value = ""
$ws.Range("E95").Value = ""
$ws.Range("F95").Value = 4701
$ws.Range("G95").Value = "新莊昌平店"
$ws.Range("H95").Value = "新北市新莊區"
$ws.Range("I95").Value = ""
$ws.Range("J95").Value = ""
$ws.Range("K95").Value = ""
$ws.Range("L95").Value = ""
$ws.Range("M95").Value = ""
$ws.Range("N95").Value = ""
$ws.Range("O95").Value = ""
$ws.Range("P95").Value = ""
$ws.Range("Q95").Value = "THILF04701"
$ws.Range("R95").Value = "新北一"
$ws.Range("S95").Value = "湯家瑋"
$ws.Range("T95").Value = 1
$ws.Range("U95").Value = "已完工"
$ws.Range("V95").Value = "2025-09-18 11:18:47"
$ws.Range("W95").Value = "2025-09-18 11:00:00"
$ws.Range("X95").Value = "2025-09-18 11:20:00"
$ws.Range("Y95").Value = ""
$ws.Range("Z95").Value = 0.3
$ws.Range("AA95").Value = ""
$ws.Range("AB95").Value = "到場處理"
$ws.Range("AC95").Value = "PMQ3+TVV"
$ws.Range("AD95").Value = "O"
$ws.Range("AE95").Value = ""
$ws.Range("AF95").Value = ""
$ws.Range("AG95").Value = ""
$ws.Range("AH95").Value = ""
$ws.Range("AI95").Value = ""
$ws.Range("AJ95").Value = "O"
$ws.Range("AK95").Value = "O"

$ws.Range("A96").Value = 94
$ws.Range("B96").Value = "服務"
$ws.Range("C96").Value = 2025092404
$ws.Range("D96").Value = ""
$ws.Range("E96").Value = ""
$ws.Range("F96").Value = 3358
$ws.Range("G96").Value = "新莊棒球場"
$ws.Range("H96").Value = "新北市新莊區"
$ws.Range("I96").Value = ""
$ws.Range("J96").Value = ""
$ws.Range("K96").Value = ""
$ws.Range("L96").Value = ""
$ws.Range("M96").Value = ""
$ws.Range("N96").Value = ""
$ws.Range("O96").Value = ""
$ws.Range("P96").Value = ""
$ws.Range("Q96").Value = "THILF03358"
$ws.Range("R96").Value = "新北一"
$ws.Range("S96").Value = "湯家瑋"
$ws.Range("T96").Value = 1
$ws.Range("U96").Value = "已完工"
$ws.Range("V96").Value = "2025-09-18 11:34:07"
$ws.Range("W96").Value = "2025-09-18 11:20:00"
$ws.Range("X96").Value = "2025-09-18 11:40:00"
$ws.Range("Y96").Value = ""
$ws.Range("Z96").Value = 0.3
$ws.Range("AA96").Value = ""
$ws.Range("AB96").Value = "到場處理"
$ws.Range("AC96").Value = "PMQ3"
$ws.Range("AD96").Value = "O"
$ws.Range("AE96").Value = ""
$ws.Range("AF96").Value = ""
$ws.Range("AG96").Value = ""
$ws.Range("AH96").Value = ""
$ws.Range("AI96").Value = ""
$ws.Range("AJ96").Value = ""
$ws.Range("AK96").Value = "O"

$ws.Range("A97").Value = 95
$ws.Range("B97").Value = "服務"
$ws.Range("C97").Value = 2025092409
$ws.Range("D97").Value = ""
$ws.Range("E97").Value = ""
$ws.Range("F97").Value = 3851
$ws.Range("G97").Value = "新莊輔園店"
$ws.Range("H97").Value = "新北市新莊區"
$ws.Range("I97").Value = ""
$ws.Range("J97").Value = ""
$ws.Range("K97").Value = ""
$ws.Range("L97").Value = ""
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = ""
$ws.Range("O97").Value = ""
$ws.Range("P97").Value = ""
$ws.Range("Q97").Value = "THILF03851"
$ws.Range("R97").Value = "新北一"
$ws.Range("S97").Value = "湯家瑋"
$ws.Range("T97").Value = 1
$ws.Range("U97").Value = "已完工"
$ws.Range("V97").Value = "2025-09-18 13:11:21"
$ws.Range("W97").Value = "2025-09-18 13:00:00"
$ws.Range("X97").Value = "2025-09-18 13:20:00"
$ws.Range("Y97").Value = ""
$ws.Range("Z97").Value = 0.3
$ws.Range("AA97").Value = ""
$ws.Range("AB97").Value = "到場處理"
$ws.Range("AC97").Value = "PMQ3+TVV"
$ws.Range("AD97").Value = "O"
$ws.Range("AE97").Value = ""
$ws.Range("AF97").Value = ""
$ws.Range("AG97").Value = ""
$ws.Range("AH97").Value = ""
$ws.Range("AI97").Value = ""
$ws.Range("AJ97").Value = "O"
$ws.Range("AK97").Value = "O"

# ---------------------------------------------------------------------
# 4. Restore wrap-text on the long-text columns (P, AC) for the rows
#    that need it (row 97's P/AC stay non-wrapping, matching the
#    template it was copied from).
# ---------------------------------------------------------------------
$ws.Range("P93").WrapText = $true
$ws.Range("AC93").WrapText = $true
$ws.Range("P94").WrapText = $true
$ws.Range("AC94").WrapText = $true
$ws.Range("P95").WrapText = $true
$ws.Range("AC95").WrapText = $true
$ws.Range("P96").WrapText = $true
$ws.Range("AC96").WrapText = $true

# ---------------------------------------------------------------------
# 5. Extend the print area to the new last row and move the saved
#    selection to A97 (matching the author's cursor position after
#    entering the new rows).
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = "'Report'!`$A`$1:`$AK`$97"
$ws.Range("A97").Select()
